$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.991.74"
$ws.Range("E2").Value = "  +0.98%  "

$ws.Range("D3").Value = "1.650.60"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").Value = "'218.23"
$ws.Range("E5").Value = "  +0.53%  "

$ws.Range("E6").Value = "  +0.56%  "

$ws.Range("E7").Value = "  -0.07%  "

$ws.Range("E8").Value = "  +1.72%  "

$ws.Range("D9").Value = "'0.0622"
$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").Value = "'19.82"
$ws.Range("E10").Value = "  +4.18%  "

$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("D12").Value = "1.884.34"
$ws.Range("E12").Value = "  +1.05%  "

$ws.Range("D13").Value = "1.655.11"
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("D14").Value = "'4.15"
$ws.Range("E14").Value = "  +0.39%  "

$ws.Range("E15").Value = "  +0.90%  "

$ws.Range("D16").Value = "'66.76"
$ws.Range("E16").Value = "  +3.71%  "

$ws.Range("D17").Value = "26.974.55"
$ws.Range("E17").Value = "  +1.01%  "

$ws.Range("D18").Value = "0.0₃0732"

$ws.Range("D19").Value = "'220.49"
$ws.Range("E19").Value = "  +4.77%  "

$ws.Range("D20").Value = "'1.01"
$ws.Range("E20").Value = "  -0.02%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'4.41"
$ws.Range("E21").Value = "  +2.00%  "

$ws.Range("B22").Value = "Chainlink"
$ws.Range("C22").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D22").Value = "'6.68"
$ws.Range("E22").Value = "  +8.13%  "

$ws.Range("E23").Value = "  +3.11%  "

$ws.Range("E24").Value = "  -0.46%  "

$ws.Range("D25").Value = "'146.42"
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").Value = "'7.38"
$ws.Range("E27").Value = "  +4.69%  "

$ws.Range("E28").Value = "  +1.12%  "

$ws.Range("E29").Value = "  +2.64%  "

$ws.Range("D30").Value = "'0.0512"
$ws.Range("E30").Value = "  +1.73%  "

$ws.Range("E31").Value = "  +0.73%  "

$ws.Range("E32").Value = "  +1.54%  "

$ws.Range("D33").Value = "'2.99"
$ws.Range("E33").Value = "  +0.87%  "

$ws.Range("E34").Value = "  +2.35%  "

$ws.Range("D35").Value = "'2.47"
$ws.Range("E35").Value = "  +1.28%  "

$ws.Range("D36").Value = "1.256.35"
$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("D37").Value = "'0.0176"
$ws.Range("E37").Value = "  +0.82%  "

$ws.Range("D38").Value = "'0.533"
$ws.Range("E38").Value = "  +0.74%  "

$ws.Range("D39").Value = "'0.834"
$ws.Range("E39").Value = "  +3.31%  "

$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("D41").Value = "'0.814"
$ws.Range("E41").Value = "  +1.62%  "

$ws.Range("D42").Value = "'5.37"
$ws.Range("E42").Value = "  +2.18%  "

$ws.Range("D43").Value = "1.795.90"
$ws.Range("E43").Value = "  +1.25%  "

$ws.Range("D44").Value = "'2.10"
$ws.Range("E44").Value = "  -4.56%  "

$ws.Range("D45").Value = "'61.43"
$ws.Range("E45").Value = "  +1.54%  "

$ws.Range("D46").Value = "'91.67"
$ws.Range("E46").Value = "  +0.71%  "

$ws.Range("E47").Value = "  +1.75%  "

$ws.Range("D48").Value = "'0.0515"
$ws.Range("E48").Value = "  -0.87%  "

$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D49").Value = "0.0₇0980"
$ws.Range("E49").Value = "  -5.78%  "

$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0976"
$ws.Range("E50").Value = "  +1.81%  "

$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.61"
$ws.Range("E51").Value = "  +1.22%  "
